# Add one new weekly price record to the "Brócoli" dataset.
#
# The sheet is a flat table (row 1 = headers, rows 2..437 = data) kept in
# reverse-chronological-ish order. A new record is inserted as the new
# row 361, pushing the former rows 361..437 down to 362..438 (which is why
# the sheet's used range grows from A1:R437 to A1:R438).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift rows 361..437 down to 362..438, leaving a blank row at 361.
$ws.Rows.Item(361).Insert()

# Populate the newly inserted row 361 with the new data point.
$ws.Range("A361").Value = 4
$ws.Range("B361").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C361").Value = "Los Lagos"
$ws.Range("D361").Value = 44889
$ws.Range("E361").Value = 10
$ws.Range("F361").Value = 100112023
$ws.Range("G361").Value = "Brócoli"
$ws.Range("H361").Value = "Sin especificar"
$ws.Range("I361").Value = "Primera"
$ws.Range("J361").Value = 600
$ws.Range("K361").Value = 1500
$ws.Range("L361").Value = 1500
$ws.Range("M361").Value = 1500
$ws.Range("N361").Value = "$/unidad"
$ws.Range("O361").Value = "Región Metropolitana"
$ws.Range("P361").Value = 1500
$ws.Range("Q361").Value = 1
$ws.Range("R361").Value = "Hortaliza"
